# "Se arreglo la validacion de modificar productos."
# Fix/update rows in the "Productos" sheet: edit a few existing product rows
# and append three newly-registered products.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Productos")

# --- widen a couple of columns (Referencia / Marca) ---------------------
# ColumnWidth includes ~0.8333 chars of cell padding on top of the raw
# character width stored in the XML, so back that padding out to land on
# an exact width of 13 / 11.
$ws.Columns.Item(1).ColumnWidth = 12.166666666666666
$ws.Columns.Item(3).ColumnWidth = 10.166666666666666

# --- correct a handful of existing product records -----------------------
$ws.Cells.Item(2, 3).Value = "sdfsdfs"
$ws.Cells.Item(2, 4).Value = 4234
$ws.Cells.Item(2, 5).Value = 2342

$ws.Cells.Item(11, 3).Value = "fsdffsf"
$ws.Cells.Item(11, 4).Value = 234234
$ws.Cells.Item(11, 5).Value = 4324234

$ws.Cells.Item(13, 3).Value = "Funciona"
$ws.Cells.Item(13, 4).Value = 11111
$ws.Cells.Item(13, 5).Value = 33333

# --- append newly registered products -------------------------------------
$ws.Cells.Item(14, 1).Value = "El Circulo"
$ws.Cells.Item(14, 2).Value = 1112232121312
$ws.Cells.Item(14, 3).Value = "Dog"
$ws.Cells.Item(14, 4).Value = 222222
$ws.Cells.Item(14, 5).Value = 111111
$ws.Cells.Item(14, 6).Value = 1221
$ws.Cells.Item(14, 7).Value = $true
$ws.Cells.Item(14, 8).Value = "04/06/2024 06:54"

$ws.Cells.Item(15, 1).Value = "The Chronic"
$ws.Cells.Item(15, 2).Value = 1123213432133
$ws.Cells.Item(15, 3).Value = "Aftermath"
$ws.Cells.Item(15, 4).Value = 100000
$ws.Cells.Item(15, 5).Value = 10000000
$ws.Cells.Item(15, 6).Value = 3333
$ws.Cells.Item(15, 7).Value = $true
$ws.Cells.Item(15, 8).Value = "04/06/2024 07:01"

$ws.Cells.Item(16, 1).Value = "Perchalsts"
# Barcode too long to stay numeric-safe, so it was entered as text
# (leading apostrophe, same as typing it directly into the cell).
$ws.Cells.Item(16, 2).Value = "'3423423432432"
$ws.Cells.Item(16, 3).Value = "sdffsdf"
$ws.Cells.Item(16, 4).Value = 4234234
$ws.Cells.Item(16, 5).Value = 234234
$ws.Cells.Item(16, 6).Value = 234234
$ws.Cells.Item(16, 7).Value = $true
$ws.Cells.Item(16, 8).Value = "04/06/2024 07:16"
